$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap rows 53 (Huelva) and 54 (Huesca) - both the province name and the
# "Casos activos" figure need to move together so the data stays attached
# to the correct province after reordering alphabetically.
$rowA = $ws.Range("A53:E53").Value()
$rowB = $ws.Range("A54:E54").Value()

$ws.Range("A53").Value = $rowB[1,1]
$ws.Range("B53").Value = $rowB[1,2]
$ws.Range("C53").Value = $rowB[1,3]
$ws.Range("D53").Value = $rowB[1,4]
$ws.Range("E53").Value = $rowB[1,5]

$ws.Range("A54").Value = $rowA[1,1]
$ws.Range("B54").Value = $rowA[1,2]
$ws.Range("C54").Value = $rowA[1,3]
$ws.Range("D54").Value = $rowA[1,4]
$ws.Range("E54").Value = $rowA[1,5]

# Update the "last updated" timestamp string (cell A1) from 09:16 to 09:46
$ws.Range("A1").Value = "Datos actualizados a 22 de Marzo de 2020 a las 09:46"
